$wb = $excel.ActiveWorkbook

# The F (and occasionally G) column values were refreshed (new scraped stats)
# on both the "展览" sheet and the aggregated "全部类型" sheet, rows 2-21.
$updates = @(
    @{ Row = 2;  F = 388 },
    @{ Row = 3;  F = 1015 },
    @{ Row = 4;  F = 253 },
    @{ Row = 5;  F = 1384; G = 58 },
    @{ Row = 6;  F = 8461 },
    @{ Row = 7;  F = 55 },
    @{ Row = 8;  F = 476 },
    @{ Row = 9;  F = 624 },
    @{ Row = 10; F = 238 },
    @{ Row = 12; F = 3407 },
    @{ Row = 13; F = 46 },
    @{ Row = 14; F = 340 },
    @{ Row = 15; F = 55 },
    @{ Row = 16; F = 938 },
    @{ Row = 17; F = 140 },
    @{ Row = 18; F = 1092 },
    @{ Row = 20; F = 156 },
    @{ Row = 21; F = 2025 }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $ws.Cells.Item($u.Row, 6).Value = $u.F
        if ($u.ContainsKey("G")) {
            $ws.Cells.Item($u.Row, 7).Value = $u.G
        }
    }
}
